# Update the SplashKit C# code samples to SplashKit C (splashkit.h) equivalents
# across every slide's "while loop" code textbox (Group 11 > TextBox 15).

$p = $ppt.ActivePresentation

$replacements = @(
    @("using static ", ""),
    @("SplashKitSDK.SplashKit", '#include "splashkit.h"'),
    @("OpenWindow", "open_window"),
    @("ClearScreen", "clear_screen"),
    @("ColorYellow", "color_yellow"),
    @("QuitRequested", "quit_requested"),
    @("FillCircle", "fill_circle"),
    @("RandomColor", "random_color"),
    @("Rnd", "rnd"),
    @("ScreenWidth", "screen_width"),
    @("ScreenHeight", "screen_height"),
    @("RefreshScreen", "refresh_screen"),
    @("ProcessEvents", "process_events")
)

for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)
    $group = $slide.Shapes.Item("Group 11")
    $textBox = $group.GroupItems.Item("TextBox 15")
    $textRange = $textBox.TextFrame.TextRange

    foreach ($pair in $replacements) {
        $textRange.Replace($pair[0], $pair[1]) | Out-Null
    }
}
